$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly-played match rows appended to the veto/pick log (rows 368-389).
# Row 368
$ws.Range("A368").Value = 27322
$ws.Range("B368").Value = 45444
$ws.Range("C368").Value = 4
$ws.Range("D368").Value = "Los Angeles Guerrillas"
$ws.Range("E368").Value = "6 Star"
$ws.Range("F368").Value = "Karachi"
$ws.Range("G368").Value = "Karachi"
$ws.Range("H368").Value = "Invasion"
$ws.Range("I368").Value = "Karachi"

# Row 369
$ws.Range("A369").Value = 27322
$ws.Range("B369").Value = 45444
$ws.Range("C369").Value = 4
$ws.Range("D369").Value = "Los Angeles Thieves"
$ws.Range("E369").Value = "Sub Base"
$ws.Range("F369").Value = "Rio"
$ws.Range("G369").Value = "6 Star"
$ws.Range("H369").Value = "Rio"
$ws.Range("J369").Value = "Invasion"

# Row 370
$ws.Range("A370").Value = 27323
$ws.Range("B370").Value = 45444
$ws.Range("C370").Value = 4
$ws.Range("D370").Value = "Las Vegas Legion"
$ws.Range("E370").Value = "6 Star"
$ws.Range("F370").Value = "Rio"
$ws.Range("G370").Value = "6 Star"
$ws.Range("H370").Value = "Rio"
$ws.Range("J370").Value = "Highrise"

# Row 371
$ws.Range("A371").Value = 27323
$ws.Range("B371").Value = 45444
$ws.Range("C371").Value = 4
$ws.Range("D371").Value = "Toronto Ultra"
$ws.Range("E371").Value = "Vista"
$ws.Range("F371").Value = "Karachi"
$ws.Range("G371").Value = "Highrise"
$ws.Range("H371").Value = "Karachi"
$ws.Range("I371").Value = "Invasion"

# Row 372
$ws.Range("A372").Value = 27324
$ws.Range("B372").Value = 45444
$ws.Range("C372").Value = 4
$ws.Range("D372").Value = "Miami Heretics"
$ws.Range("E372").Value = "Vista"
$ws.Range("F372").Value = "Karachi"
$ws.Range("G372").Value = "Rio"
$ws.Range("H372").Value = "Highrise"
$ws.Range("J372").Value = "Invasion"

# Row 373
$ws.Range("A373").Value = 27324
$ws.Range("B373").Value = 45444
$ws.Range("C373").Value = 4
$ws.Range("D373").Value = "OpTic Texas"
$ws.Range("E373").Value = "6 Star"
$ws.Range("F373").Value = "Sub Base"
$ws.Range("G373").Value = "6 Star"
$ws.Range("H373").Value = "Karachi"
$ws.Range("I373").Value = "Karachi"

# Row 374
$ws.Range("A374").Value = 27325
$ws.Range("B374").Value = 45444
$ws.Range("C374").Value = 4
$ws.Range("D374").Value = "Carolina Royal Ravens"
$ws.Range("E374").Value = "Rio"
$ws.Range("F374").Value = "Vista"
$ws.Range("G374").Value = "Highrise"
$ws.Range("H374").Value = "Rio"
$ws.Range("J374").Value = "Invasion"

# Row 375
$ws.Range("A375").Value = 27325
$ws.Range("B375").Value = 45444
$ws.Range("C375").Value = 4
$ws.Range("D375").Value = "New York Subliners"
$ws.Range("E375").Value = "Sub Base"
$ws.Range("F375").Value = "Karachi"
$ws.Range("G375").Value = "Karachi"
$ws.Range("H375").Value = "Invasion"
$ws.Range("I375").Value = "Karachi"

# Row 376
$ws.Range("A376").Value = 27326
$ws.Range("B376").Value = 45445
$ws.Range("C376").Value = 4
$ws.Range("D376").Value = "Minnesota ROKKR"
$ws.Range("E376").Value = "Vista"
$ws.Range("F376").Value = "Rio"
$ws.Range("G376").Value = "Invasion"
$ws.Range("H376").Value = "Rio"
$ws.Range("J376").Value = "Invasion"

# Row 377
$ws.Range("A377").Value = 27326
$ws.Range("B377").Value = 45445
$ws.Range("C377").Value = 4
$ws.Range("D377").Value = "Seattle Surge"
$ws.Range("E377").Value = "Sub Base"
$ws.Range("F377").Value = "6 Star"
$ws.Range("G377").Value = "Highrise"
$ws.Range("H377").Value = "6 Star"
$ws.Range("I377").Value = "Karachi"

# Row 378
$ws.Range("A378").Value = 27327
$ws.Range("B378").Value = 45445
$ws.Range("C378").Value = 4
$ws.Range("D378").Value = "Boston Breach"
$ws.Range("E378").Value = "6 Star"
$ws.Range("F378").Value = "Vista"
$ws.Range("G378").Value = "Highrise"
$ws.Range("H378").Value = "Invasion"
$ws.Range("I378").Value = "Invasion"

# Row 379
$ws.Range("A379").Value = 27327
$ws.Range("B379").Value = 45445
$ws.Range("C379").Value = 4
$ws.Range("D379").Value = "Los Angeles Thieves"
$ws.Range("E379").Value = "Karachi"
$ws.Range("F379").Value = "Rio"
$ws.Range("G379").Value = "6 Star"
$ws.Range("H379").Value = "Rio"
$ws.Range("J379").Value = "Highrise"

# Row 380
$ws.Range("A380").Value = 27328
$ws.Range("B380").Value = 45445
$ws.Range("C380").Value = 4
$ws.Range("D380").Value = "Atlanta FaZe"
$ws.Range("E380").Value = "Sub Base"
$ws.Range("F380").Value = "Vista"
$ws.Range("G380").Value = "Highrise"
$ws.Range("H380").Value = "6 Star"
$ws.Range("I380").Value = "Highrise"

# Row 381
$ws.Range("A381").Value = 27328
$ws.Range("B381").Value = 45445
$ws.Range("C381").Value = 4
$ws.Range("D381").Value = "New York Subliners"
$ws.Range("E381").Value = "Rio"
$ws.Range("F381").Value = "Karachi"
$ws.Range("G381").Value = "Rio"
$ws.Range("H381").Value = "Invasion"
$ws.Range("J381").Value = "Karachi"

# Row 382
$ws.Range("A382").Value = 27329
$ws.Range("B382").Value = 45445
$ws.Range("C382").Value = 4
$ws.Range("D382").Value = "Las Vegas Legion"
$ws.Range("E382").Value = "6 Star"
$ws.Range("F382").Value = "Vista"
$ws.Range("G382").Value = "Invasion"
$ws.Range("H382").Value = "6 Star"
$ws.Range("J382").Value = "Highrise"

# Row 383
$ws.Range("A383").Value = 27329
$ws.Range("B383").Value = 45445
$ws.Range("C383").Value = 4
$ws.Range("D383").Value = "OpTic Texas"
$ws.Range("E383").Value = "Karachi"
$ws.Range("F383").Value = "Sub Base"
$ws.Range("G383").Value = "Rio"
$ws.Range("H383").Value = "Karachi"
$ws.Range("I383").Value = "Invasion"

# Row 384
$ws.Range("A384").Value = 27330
$ws.Range("B384").Value = 45450
$ws.Range("C384").Value = 4
$ws.Range("D384").Value = "Seattle Surge"
$ws.Range("E384").Value = "Karachi"
$ws.Range("F384").Value = "Vista"
$ws.Range("G384").Value = "Highrise"
$ws.Range("H384").Value = "6 Star"
$ws.Range("I384").Value = "Karachi"

# Row 385
$ws.Range("A385").Value = 27330
$ws.Range("B385").Value = 45450
$ws.Range("C385").Value = 4
$ws.Range("D385").Value = "New York Subliners"
$ws.Range("E385").Value = "6 Star"
$ws.Range("F385").Value = "Rio"
$ws.Range("G385").Value = "Rio"
$ws.Range("H385").Value = "Invasion"
$ws.Range("J385").Value = "Highrise"

# Row 386
$ws.Range("A386").Value = 27331
$ws.Range("B386").Value = 45450
$ws.Range("C386").Value = 4
$ws.Range("D386").Value = "Miami Heretics"
$ws.Range("E386").Value = "Rio"
$ws.Range("F386").Value = "6 Star"
$ws.Range("G386").Value = "Rio"
$ws.Range("H386").Value = "Highrise"
$ws.Range("I386").Value = "Invasion"

# Row 387
$ws.Range("A387").Value = 27331
$ws.Range("B387").Value = 45450
$ws.Range("C387").Value = 4
$ws.Range("D387").Value = "Minnesota ROKKR"
$ws.Range("E387").Value = "Sub Base"
$ws.Range("F387").Value = "Vista"
$ws.Range("G387").Value = "6 Star"
$ws.Range("H387").Value = "Karachi"
$ws.Range("J387").Value = "Highrise"

# Row 388
$ws.Range("A388").Value = 27332
$ws.Range("B388").Value = 45450
$ws.Range("C388").Value = 4
$ws.Range("D388").Value = "Carolina Royal Ravens"
$ws.Range("E388").Value = "Vista"
$ws.Range("F388").Value = "Sub Base"
$ws.Range("G388").Value = "6 Star"
$ws.Range("H388").Value = "Rio"
$ws.Range("I388").Value = "Highrise"

# Row 389
$ws.Range("A389").Value = 27332
$ws.Range("B389").Value = 45450
$ws.Range("C389").Value = 4
$ws.Range("D389").Value = "Atlanta FaZe"
$ws.Range("E389").Value = "Karachi"
$ws.Range("F389").Value = "6 Star"
$ws.Range("G389").Value = "Karachi"
$ws.Range("H389").Value = "Invasion"
$ws.Range("J389").Value = "Invasion"

# Restore the last-used cell selection after entering the new rows.
$ws.Range("L392").Select()
